# NPs understand byte array messages. Full timings recorded
#
# Adds a new "NAMED PIPE RESULTS" worksheet (mirroring the layout of the
# existing "MAILSLOT RESULTS" sheet) with the full named-pipe timing data,
# and updates the selection/active-sheet state accordingly.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("MAILSLOT RESULTS")

# --- 1. Add the new worksheet at the end of the workbook -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "NAMED PIPE RESULTS"

# --- 2. Clone the look of the "MAILSLOT RESULTS" sheet's first table -------
# Title band (merged, bold/underline, yellow fill, centered)
$ws2.Range("A1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)

# Column header row (bold, centered)
$ws2.Range("B2:E2").Copy()
$ws.Range("B2:E2").PasteSpecial(-4122)

# Row-number column (left aligned)
$ws2.Range("A3:A12").Copy()
$ws.Range("A3:A12").PasteSpecial(-4122)

# "Average" summary row (bold)
$ws2.Range("A13:E13").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 3. Title / headers ------------------------------------------------
$ws.Range("A1").Value = "NAMED PIPES "
$ws.Range("A1:E1").Merge()

$ws.Range("B2").Value = "40 BYTES "
$ws.Range("C2").Value = "400 BYTES"
$ws.Range("D2").Value = "4000 BYTES"
$ws.Range("E2").Value = "40 000 BYTES"

# --- 4. Timing data (10 runs x 4 message sizes) -----------------------------
$bVals = @(672093, 773663, 995681, 1108577, 1071574, 680401, 1064023, 1051940, 1041746, 618855)
$cVals = @(818218, 1051562, 1053073, 926206, 842383, 1053828, 761203, 845781, 1148600, 941687)
$dVals = @(890335, 870324, 835964, 797073, 1341544, 1311338, 714760, 805002, 593934, 807645)
$eVals = @(1057604, 958300, 637734, 1122925, 723067, 998324, 1127834, 968117, 1051562, 1038347)

for ($i = 0; $i -lt 10; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}

# --- 5. Average summary row --------------------------------------------
$ws.Range("A13").Value = "Average"
$ws.Range("B13").Formula = "=AVERAGE(B3:B12)"
$ws.Range("C13").Formula = "=AVERAGE(C3:C12)"
$ws.Range("D13").Formula = "=AVERAGE(D3:D12)"
$ws.Range("E13").Formula = "=AVERAGE(E3:E12)"

# --- 6. Column width for the "40 000 BYTES" column --------------------
$ws.Columns.Item(5).ColumnWidth = 10.39

# --- 7. Selection / active-sheet bookkeeping --------------------------
# "MAILSLOT RESULTS" keeps a fresh selection but is no longer the active tab
$ws2.Activate()
$ws2.Range("A3:A13").Select()

# The new "NAMED PIPE RESULTS" sheet becomes the active / selected tab
$ws.Activate()
$ws.Range("F16").Select()
